# Re-upload of the daily "Saldo" client-balance export: the reference date
# rolled forward by one day (2024-10-29 -> 2024-10-30) and a handful of
# client balances were refreshed with newer figures pulled for that date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is named after the export run (date + time) -> rename it to
# match the new export timestamp.
$ws.Name = "IClientBalance-20241030-105238-"

# Column G ("Dt. Referencia") holds the reference date as a serial number.
# Every data row (2..274) moves from 45594 (2024-10-29) to 45595 (2024-10-30).
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45595
}

# A subset of rows also got refreshed "Saldo Previsto" (E) / "Vl. Total" (H)
# amounts (and, for row 224, a refreshed "Vl. Projetado" (D) amount too) as
# part of this re-export. Row -> column letter -> new value.
$updates = @(
    @{ Row = 6;   Col = 5; Value = 1105.22 },
    @{ Row = 6;   Col = 8; Value = 1105.22 },

    @{ Row = 15;  Col = 5; Value = 971.57 },
    @{ Row = 15;  Col = 8; Value = 971.57 },

    @{ Row = 17;  Col = 5; Value = 1152.19 },
    @{ Row = 17;  Col = 8; Value = 1152.19 },

    @{ Row = 43;  Col = 5; Value = 1456.6 },
    @{ Row = 43;  Col = 8; Value = 1456.6 },

    @{ Row = 49;  Col = 5; Value = 1304.4000000000001 },
    @{ Row = 49;  Col = 8; Value = 1304.4000000000001 },

    @{ Row = 51;  Col = 5; Value = 740.68 },
    @{ Row = 51;  Col = 8; Value = 740.68 },

    @{ Row = 52;  Col = 5; Value = 1336.19 },
    @{ Row = 52;  Col = 8; Value = 1336.19 },

    @{ Row = 57;  Col = 5; Value = 2422.89 },
    @{ Row = 57;  Col = 8; Value = 2422.89 },

    @{ Row = 97;  Col = 5; Value = 76.099999999999994 },
    @{ Row = 97;  Col = 8; Value = 76.099999999999994 },

    @{ Row = 101; Col = 5; Value = 152.19999999999999 },
    @{ Row = 101; Col = 8; Value = 152.19999999999999 },

    @{ Row = 102; Col = 5; Value = 1228.3 },
    @{ Row = 102; Col = 8; Value = 1228.3 },

    @{ Row = 105; Col = 5; Value = 845.81 },
    @{ Row = 105; Col = 8; Value = 845.81 },

    @{ Row = 107; Col = 5; Value = 33258.089999999997 },
    @{ Row = 107; Col = 8; Value = 33258.089999999997 },

    @{ Row = 108; Col = 5; Value = 0 },
    @{ Row = 108; Col = 8; Value = 0 },

    @{ Row = 109; Col = 5; Value = 0 },
    @{ Row = 109; Col = 8; Value = 0 },

    @{ Row = 110; Col = 5; Value = 897.16 },
    @{ Row = 110; Col = 8; Value = 897.16 },

    @{ Row = 112; Col = 5; Value = 2.38 },
    @{ Row = 112; Col = 8; Value = 2.38 },

    @{ Row = 120; Col = 5; Value = 997.03 },
    @{ Row = 120; Col = 8; Value = 997.03 },

    @{ Row = 138; Col = 5; Value = 1861.8 },
    @{ Row = 138; Col = 8; Value = 1861.8 },

    @{ Row = 143; Col = 5; Value = 13470.71 },
    @{ Row = 143; Col = 8; Value = 13470.71 },

    # Row 224 is the odd one out: it's "Vl. Projetado" (D) that changed
    # (0 -> 9915.9), "Saldo Previsto" (E) stays at 643.46, and "Vl. Total"
    # (H) reflects the new sum (9915.9 + 643.46 = 10559.36).
    @{ Row = 224; Col = 4; Value = 9915.9 },
    @{ Row = 224; Col = 8; Value = 10559.36 },

    @{ Row = 230; Col = 5; Value = 843.5 },
    @{ Row = 230; Col = 8; Value = 843.5 },

    @{ Row = 255; Col = 5; Value = 1684.91 },
    @{ Row = 255; Col = 8; Value = 1684.91 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}
